# Update gh-pages to output generated at 456a3b4
# Applies numeric updates to column F ("views"/count column) across sheets.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F10").Value = 6778
$ws1.Range("F11").Value = 31
$ws1.Range("F13").Value = 328
$ws1.Range("F21").Value = 93
$ws1.Range("F22").Value = 1089
$ws1.Range("F25").Value = 318
$ws1.Range("F26").Value = 1649
$ws1.Range("F27").Value = 1634
$ws1.Range("F29").Value = 699
$ws1.Range("F35").Value = 18
$ws1.Range("F43").Value = 13

# --- Sheet "演出" ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F20").Value = 17
$ws2.Range("F23").Value = 452

# --- Sheet "本地生活" ---
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F14").Value = 1268
$ws3.Range("F15").Value = 7042

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F12").Value = 6778
$ws4.Range("F15").Value = 31
$ws4.Range("F16").Value = 328
$ws4.Range("F23").Value = 93
$ws4.Range("F24").Value = 1089
$ws4.Range("F26").Value = 318
$ws4.Range("F28").Value = 1649
$ws4.Range("F35").Value = 452
